$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change, label stays the same
$ws.Range("B3").Value = 0.9551407937152788
$ws.Range("C3").Value = 0.9541076839629467
$ws.Range("D3").Value = 0.935751436688177

# Row 4: label changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.8885420669285531
$ws.Range("C4").Value = 0.891202328577284
$ws.Range("D4").Value = 0.6617280785679465

# Row 5: label changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9143669643010739
$ws.Range("C5").Value = 0.9121972328564959
$ws.Range("D5").Value = 0.9042129056918151
